$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 155.57143
$ws.Range("I33").Value = 145.23077
$ws.Range("K33").Value = 145.23077
$ws.Range("M33").Value = 83.76922999999999

$ws.Range("H58").Value = 1865
$ws.Range("I58").Value = 514.25
$ws.Range("J58").Value = 3666
$ws.Range("K58").Value = 1542.75
$ws.Range("L58").Value = 10998
$ws.Range("M58").Value = -1392.75
$ws.Range("N58").Value = -11298

$ws.Range("H87").Value = 27439.96
$ws.Range("J87").Value = 27439.96
$ws.Range("L87").Value = 27439.96
$ws.Range("N87").Value = -29935.96

$ws.Range("H90").Value = 27439.96
$ws.Range("J90").Value = 27439.96
$ws.Range("L90").Value = 82319.88
$ws.Range("N90").Value = -94799.88

$ws.Range("H112").Value = 1613.2222
$ws.Range("J112").Value = 1710.0714
$ws.Range("L112").Value = 5130.2142
$ws.Range("N112").Value = -7346.2142

$ws.Range("H130").Value = 42500
$ws.Range("J130").Value = 42500
$ws.Range("L130").Value = 42500
$ws.Range("N130").Value = -52540

$ws.Range("H132").Value = 9417.183999999999
$ws.Range("I132").Value = 986.5088
$ws.Range("J132").Value = 169600
$ws.Range("K132").Value = 2959.5264
$ws.Range("L132").Value = 508800
$ws.Range("M132").Value = -429.5263999999997
$ws.Range("N132").Value = -513860

$ws.Range("H138").Value = 4712.905
$ws.Range("I138").Value = 2999
$ws.Range("J138").Value = 4998.5557
$ws.Range("K138").Value = 8997
$ws.Range("L138").Value = 14995.6671
$ws.Range("M138").Value = -3857
$ws.Range("N138").Value = -25275.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13030.866
$ws.Range("I2").Value = 15327.958
$ws.Range("K2").Value = 15327.958
$ws.Range("M2").Value = -15214.958

$ws.Range("H34").Value = 86631.14
$ws.Range("I34").Value = 12847.5
$ws.Range("J34").Value = 185009.33
$ws.Range("K34").Value = 12847.5
$ws.Range("L34").Value = 185009.33
$ws.Range("M34").Value = -12576.5
$ws.Range("N34").Value = -185551.33

$ws.Range("H44").Value = 44747.5
$ws.Range("J44").Value = 44747.5
$ws.Range("L44").Value = 44747.5
$ws.Range("N44").Value = -45723.5

$ws.Range("H55").Value = 181666.56
$ws.Range("I55").Value = 34999
$ws.Range("J55").Value = 200000
$ws.Range("K55").Value = 34999
$ws.Range("L55").Value = 200000
$ws.Range("M55").Value = -34684
$ws.Range("N55").Value = -200630

$ws.Range("H61").Value = 3320.7646
$ws.Range("I61").Value = 3748.3
$ws.Range("K61").Value = 3748.3
$ws.Range("M61").Value = -3536.3

$ws.Range("H74").Value = 1680.95
$ws.Range("I74").Value = 1564.0625
$ws.Range("K74").Value = 1564.0625
$ws.Range("M74").Value = -690.0625

$ws.Range("H77").Value = 1680.95
$ws.Range("I77").Value = 1564.0625
$ws.Range("K77").Value = 7820.3125
$ws.Range("M77").Value = -3452.3125

$ws.Range("H110").Value = 3284.0425
$ws.Range("I110").Value = 3089.0938
$ws.Range("K110").Value = 3089.0938
$ws.Range("M110").Value = -1044.0938

$ws.Range("H116").Value = 13030.866
$ws.Range("I116").Value = 15327.958
$ws.Range("K116").Value = 15327.958
$ws.Range("M116").Value = -13033.958

$ws.Range("H136").Value = 3320.7646
$ws.Range("I136").Value = 3748.3
$ws.Range("K136").Value = 11244.9
$ws.Range("M136").Value = -8694.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13030.866
$ws.Range("I3").Value = 15327.958
$ws.Range("K3").Value = 15327.958
$ws.Range("M3").Value = -15213.958

$ws.Range("H80").Value = 352.13333
$ws.Range("J80").Value = 369.3846
$ws.Range("L80").Value = 369.3846
$ws.Range("N80").Value = -2365.3846

$ws.Range("H83").Value = 352.13333
$ws.Range("J83").Value = 369.3846
$ws.Range("L83").Value = 1846.923
$ws.Range("N83").Value = -11830.923

$ws.Range("H86").Value = 7469.4707
$ws.Range("I86").Value = 6756.143
$ws.Range("J86").Value = 10798.333
$ws.Range("K86").Value = 6756.143
$ws.Range("L86").Value = 10798.333
$ws.Range("M86").Value = -5633.143
$ws.Range("N86").Value = -13044.333

$ws.Range("H89").Value = 7469.4707
$ws.Range("I89").Value = 6756.143
$ws.Range("J89").Value = 10798.333
$ws.Range("K89").Value = 33780.715
$ws.Range("L89").Value = 53991.665
$ws.Range("M89").Value = -28164.715
$ws.Range("N89").Value = -65223.665

$ws.Range("H94").Value = 2425.4707
$ws.Range("J94").Value = 2912.25
$ws.Range("L94").Value = 2912.25
$ws.Range("N94").Value = -3814.25

$ws.Range("H105").Value = 3900.5833
$ws.Range("I105").Value = 3138.6875
$ws.Range("J105").Value = 5424.375
$ws.Range("K105").Value = 3138.6875
$ws.Range("L105").Value = 5424.375
$ws.Range("M105").Value = -1391.6875
$ws.Range("N105").Value = -8918.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 866.6667
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = -126
$ws.Range("N17").Value = -2348

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = $null

$ws.Range("H31").Value = 12904.5
$ws.Range("I31").Value = 26776
$ws.Range("K31").Value = 26776
$ws.Range("M31").Value = -26481

$ws.Range("H34").Value = 12904.5
$ws.Range("I34").Value = 26776
$ws.Range("K34").Value = 26776
$ws.Range("M34").Value = -26574

$ws.Range("H105").Value = 6749.3335
$ws.Range("J105").Value = 2999
$ws.Range("L105").Value = 2999
$ws.Range("N105").Value = -6493

$ws.Range("H107").Value = 409.42856
$ws.Range("I107").Value = 402.66666
$ws.Range("K107").Value = 402.66666
$ws.Range("M107").Value = 1517.33334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 811.3333
$ws.Range("I18").Value = 811.3333
$ws.Range("K18").Value = 2433.9999
$ws.Range("M18").Value = -2264.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5331.926
$ws.Range("I97").Value = 1638.48
$ws.Range("J97").Value = 51500
$ws.Range("K97").Value = 1638.48
$ws.Range("L97").Value = 51500
$ws.Range("M97").Value = -1142.48
$ws.Range("N97").Value = -52492

$ws.Range("H102").Value = 3028.96
$ws.Range("I102").Value = 3032.9092
$ws.Range("K102").Value = 3032.9092
$ws.Range("M102").Value = -1410.9092

$ws.Range("H132").Value = 11800
$ws.Range("I132").Value = 11000
$ws.Range("K132").Value = 33000
$ws.Range("M132").Value = -30470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2493.625
$ws.Range("I22").Value = 2292.8
$ws.Range("K22").Value = 2292.8
$ws.Range("M22").Value = -1997.8

$ws.Range("H27").Value = 2493.625
$ws.Range("I27").Value = 2292.8
$ws.Range("K27").Value = 2292.8
$ws.Range("M27").Value = -2185.8

$ws.Range("H55").Value = 524.6429000000001
$ws.Range("I55").Value = 517.8570999999999
$ws.Range("K55").Value = 517.8570999999999
$ws.Range("M55").Value = -344.8570999999999

$ws.Range("H68").Value = 3524.875
$ws.Range("I68").Value = 3524.875
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3524.875
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2775.875
$ws.Range("N68").Value = $null

$ws.Range("H71").Value = 3524.875
$ws.Range("I71").Value = 3524.875
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 17624.375
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -13880.375
$ws.Range("N71").Value = $null

$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50676

$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52340

$ws.Range("H122").Value = 2822.5
$ws.Range("I122").Value = 2763.3333
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8289.999899999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5839.999899999999
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 22191.428
$ws.Range("J41").Value = 22191.428
$ws.Range("L41").Value = 22191.428
$ws.Range("N41").Value = -22971.428

$ws.Range("H138").Value = 76750
$ws.Range("J138").Value = 76750
$ws.Range("L138").Value = 76750
$ws.Range("N138").Value = -87030
